$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.218987822532654
$ws.Range("B1").Value = 2.360851049423218
$ws.Range("C1").Value = 4.719763278961182
$ws.Range("D1").Value = 3.225871324539185
$ws.Range("E1").Value = 1.15953254699707
